$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Worksheet data updates ---
$ws.Range("A2").Value = 420
$ws.Range("B2").Value = 420
$ws.Range("C2").Value = 420
$ws.Range("D2").Value = 420

$ws.Range("D5").Value = 360
$ws.Range("E5").Value = 6

$ws.Range("A6").Value = 360
$ws.Range("C6").Value = 300
$ws.Range("D6").Value = 300
$ws.Range("E6").Value = 9

$ws.Range("G9").Value = 300

# --- Chart 1 (was "US-East-1 128MB") -> "US-East-2 2048MB", drop duplicate series, rename axes ---
$chart1 = $ws.ChartObjects(1).Chart
$chart1.ChartTitle.Text = "US-East-2 2048MB"
$chart1.SeriesCollection(2).Delete()
$chart1.Axes(1).AxisTitle.Text = "Warmstarts"
$chart1.Axes(2).AxisTitle.Text = "Seconds until Warmstart"

# --- Chart 2 (US-East-1 1024MB): rename axes only ---
$chart2 = $ws.ChartObjects(2).Chart
$chart2.Axes(1).AxisTitle.Text = "Warmstarts"
$chart2.Axes(2).AxisTitle.Text = "Seconds until Warmstart"

# --- Chart 3 (US-East-1 2048MB): no text changes ---

# --- Chart 4 (US-East-1 3008MB): rename axes only ---
$chart4 = $ws.ChartObjects(4).Chart
$chart4.Axes(1).AxisTitle.Text = "Warmstarts"
$chart4.Axes(2).AxisTitle.Text = "Seconds until Warmstart"

# --- Chart 5 (US-East-2 128MB): drop duplicate series only ---
$chart5 = $ws.ChartObjects(5).Chart
$chart5.SeriesCollection(2).Delete()

# --- Charts 6,7,8: no text changes ---

# --- Selection ---
$ws.Range("E27").Select()
